$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Compartments sheet: remove the "Can Calibrate" column (F), so that the
#    "Databook Page" column (G) shifts left into F.
# ---------------------------------------------------------------------------
$wsComp = $wb.Worksheets.Item("Compartments")
$compDatabookPageText = $wsComp.Range("G1").Comment.Text()
$wsComp.Range("F1").Comment.Delete()
$wsComp.Range("G1").Comment.Delete()
$wsComp.Range("F1").EntireColumn.Delete()
$wsComp.Range("F1").AddComment($compDatabookPageText)

# ---------------------------------------------------------------------------
# 2. Characteristics sheet: remove the "Can Calibrate" column (E), so that
#    the "Databook Page" column (F) shifts left into E.
# ---------------------------------------------------------------------------
$wsChar = $wb.Worksheets.Item("Characteristics")
$charDatabookPageText = $wsChar.Range("F1").Comment.Text()
$wsChar.Range("E1").Comment.Delete()
$wsChar.Range("F1").Comment.Delete()
$wsChar.Range("E1").EntireColumn.Delete()
$wsChar.Range("E1").AddComment($charDatabookPageText)

# ---------------------------------------------------------------------------
# 3. Parameters sheet: rename "Is Impact" (H1) to "Targetable", and remove
#    the "Can Calibrate" column (I), so that "Databook Page" (J) shifts left
#    into I.
# ---------------------------------------------------------------------------
$wsParam = $wb.Worksheets.Item("Parameters")
$wsParam.Range("H1").Value = "Targetable"
$paramDatabookPageText = $wsParam.Range("J1").Comment.Text()
$wsParam.Range("I1").Comment.Delete()
$wsParam.Range("J1").Comment.Delete()
$wsParam.Range("I1").EntireColumn.Delete()
$wsParam.Range("I1").AddComment($paramDatabookPageText)
